$d = $word.ActiveDocument

# 1. Replace " (fare la regressione " with "." in the paragraph about NCDs/regression.
$d.Content.Find.Execute(" (fare la regressione ", $true, $false, $false, $false, $false,
                         $true, 1, $false, ".", 2)

# 2. Locate the paragraph that now ends with "...NCDs." to anchor the new content after it.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "*spiega solo parzialmente la percentuale di ricoveri per NCDs.*") {
        $target = $d.Paragraphs.Item($i)
        break
    }
}
Write-Output ("target index=" + $target.Range.Start)

# 3. Insert new paragraphs after it.
$target.Range.InsertParagraphAfter()
$cur = $target.Next()
$cur.Range.InsertAfter("Come mostra la regressione c’è un associazione tra la percentuale di anziani e la percentuale di ricoveri per queste patologie, ma il modello spiega solo il X% della variabilità. Questo significa che c’è molto spazio per migliorare. ")

$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()

$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.InsertAfter("Quelli che seguono sono i comuni che si discostano di più dal valore atteso. ")

$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.InsertAfter("I primi sono quelli che a parità di anziani hanno più ricoveri, quindi dove c’è più bisogno di intervenire. ")

$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()

$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.InsertAfter("I secondi, al contrario, sono quelli che a parità di anziani generano meno ricoveri, e sarebbe interessante vedere quali potrebbero essere le cause del minor accesso all’ospedale (migliore salute? Migliori servizi?)")

$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()

$cur.Range.InsertParagraphAfter()
$cur = $cur.Next()
$cur.Range.InsertAfter("Nella mappa i primi e i secondi. ")

Write-Output ("final paragraph count=" + $d.Paragraphs.Count)
